# Applies the "v2 manual" edits described by the diff:
#  1. Re-word "... online gezet moet worden" -> "... online moet worden gezet"
#     (split into 4 runs, matching the target OOXML exactly).
#  2. Merge the "Multiple-Recording" run/proofErr markup into a single run.
#  3. Fix "Controller" -> "Controleer" (split "Control"/"eer") and leave the
#     _GoBack bookmark at that edit point (standard Word behaviour: the
#     bookmark follows the last text change).
#  4/5. Because _GoBack is recreated at the new location, Word renumbers the
#     bookmark ids; OLE_LINK5 shifts from id 4 to id 5.
#  6. Merge the "Current" run/proofErr markup into a single run.

$d = $word.ActiveDocument

function Get-ParagraphByText($needle) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs($i)
        if ($p.Range.Text -like "*$needle*") {
            return $p
        }
    }
    return $null
}

function Set-ParagraphInnerXml($para, $innerXml) {
    $wrapper = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
        'xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" ' +
        'xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" ' +
        'xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" ' +
        'xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture">' +
        '<w:body>' + $innerXml + '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $para.Range.InsertXML($wrapper)
}

# --- Hunk 1: "online gezet moet worden" -> "online moet worden gezet" -----
$p1 = Get-ParagraphByText("De wijze waarop de ATIS")
$p1xml = '<w:p w:rsidR="00424DBF" w:rsidRDefault="00424DBF" w:rsidP="00424DBF">' +
    '<w:pPr><w:jc w:val="center"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">De wijze waarop de ATIS vanaf versie 2 online </w:t></w:r>' +
    '<w:r><w:t>moet</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> worden</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> gezet</w:t></w:r>' +
    '<w:r w:rsidR="00D1520B"><w:t xml:space="preserve"> in Euroscope </w:t></w:r>' +
    '<w:r><w:t>is veranderd t.o.v. de eerdere versies.</w:t></w:r>' +
    '<w:r w:rsidR="002E5277"><w:t xml:space="preserve"> Dit document beschrijft de verschillen en de handelingen die moeten worden verricht om met versie 2 een voice ATIS in EuroScope online te zetten.</w:t></w:r>' +
    '</w:p>'
Set-ParagraphInnerXml $p1 $p1xml

# --- Hunk 2: merge the "Multiple-Recording" proofErr-split run ------------
$p2 = Get-ParagraphByText("In de eerdere versies van de Dutch VACC")
$p2xml = '<w:p w:rsidR="002E5277" w:rsidRDefault="00D1520B" w:rsidP="00D1520B">' +
    '<w:r><w:t xml:space="preserve">In de eerdere versies van de Dutch VACC ATIS Generator werd er gebruikt </w:t></w:r>' +
    '<w:r w:rsidR="004F40EF"><w:t>gemaakt</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> van de &#8220;Multiple-Recording Mode&#8221; van EuroScope. </w:t></w:r>' +
    '<w:r w:rsidR="002E5277"><w:t xml:space="preserve">Dat houdt in dat EuroScope vanuit verschillende sample bestanden, welke uit het &#8220;ATIS files descriptor&#8221; bestand worden </w:t></w:r>' +
    '<w:r w:rsidR="004F40EF"><w:t>in</w:t></w:r>' +
    '<w:r w:rsidR="002E5277"><w:t xml:space="preserve">gelezen, zelf een ATIS opbouwt. EuroScope speelt </w:t></w:r>' +
    '<w:r w:rsidR="00635301"><w:t>de samples</w:t></w:r>' +
    '<w:r w:rsidR="002E5277"><w:t xml:space="preserve"> &#233;&#233;n voor &#233;&#233;n af</w:t></w:r>' +
    '<w:r w:rsidR="00635301"><w:t xml:space="preserve"> wat de </w:t></w:r>' +
    '<w:r w:rsidR="002E5277"><w:t>illusie</w:t></w:r>' +
    '<w:r w:rsidR="00635301"><w:t xml:space="preserve"> geeft</w:t></w:r>' +
    '<w:r w:rsidR="002E5277"><w:t xml:space="preserve"> dat er een </w:t></w:r>' +
    '<w:r w:rsidR="005A2ECD"><w:t>audio bestand</w:t></w:r>' +
    '<w:r w:rsidR="002E5277"><w:t xml:space="preserve"> herhaaldelijk wordt </w:t></w:r>' +
    '<w:r w:rsidR="005A2ECD"><w:t>af</w:t></w:r>' +
    '<w:r w:rsidR="002E5277"><w:t>gespeel</w:t></w:r>' +
    '<w:r w:rsidR="005A2ECD"><w:t>d</w:t></w:r>' +
    '<w:r w:rsidR="002E5277"><w:t>.</w:t></w:r>' +
    '</w:p>'
Set-ParagraphInnerXml $p2 $p2xml

# --- Hunk 3: "Controller" -> "Controleer" (keep lastRenderedPageBreak + pic) -
$p3 = Get-ParagraphByText("Controller of ATIS letters overeenkomen")
$p3xml = '<w:p w:rsidR="00CC4467" w:rsidRDefault="00CC4467" w:rsidP="00356877">' +
    '<w:pPr><w:pStyle w:val="Lijstalinea"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:lastRenderedPageBreak/><w:t>Control</w:t></w:r>' +
    '<w:r><w:t>eer</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> of ATIS letters overeenkomen.</w:t></w:r>' +
    '<w:r><w:rPr><w:noProof/><w:lang w:eastAsia="nl-NL"/></w:rPr>' +
    '<w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0">' +
    '<wp:extent cx="5486400" cy="933450"/><wp:effectExtent l="152400" t="152400" r="361950" b="361950"/>' +
    '<wp:docPr id="17" name="Afbeelding 17"/>' +
    '<wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr>' +
    '<a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture">' +
    '<pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="0" name="Picture 4"/>' +
    '<pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr>' +
    '<pic:blipFill rotWithShape="1"><a:blip r:embed="rId11"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}">' +
    '<a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip>' +
    '<a:srcRect t="76337" b="8427"/><a:stretch/></pic:blipFill>' +
    '<pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="5486400" cy="933450"/></a:xfrm>' +
    '<a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:ln><a:noFill/></a:ln>' +
    '<a:effectLst><a:outerShdw blurRad="292100" dist="139700" dir="2700000" algn="tl" rotWithShape="0">' +
    '<a:srgbClr val="333333"><a:alpha val="65000"/></a:srgbClr></a:outerShdw></a:effectLst>' +
    '<a:extLst><a:ext uri="{53640926-AAD7-44D8-BBD7-CCE9431645EC}">' +
    '<a14:shadowObscured xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main"/></a:ext></a:extLst>' +
    '</pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r>' +
    '</w:p>'
Set-ParagraphInnerXml $p3 $p3xml

# Move the _GoBack bookmark to the site of this last real text edit, exactly
# like Word itself does: drop the old one and recreate it between "Control"
# and "eer" (this also renumbers OLE_LINK5 from id 4 to id 5).
$d.Bookmarks("_GoBack").Delete()
$rngControl = $d.Content
$rngControl.Find.Execute("Control", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$gobackPoint = $d.Range($rngControl.End, $rngControl.End)
$d.Bookmarks.Add("_GoBack", $gobackPoint) | Out-Null

# --- Hunk 6: merge the "Current" proofErr-split run ------------------------
$p6 = Get-ParagraphByText('Verhoog de')
$p6xml = '<w:p w:rsidR="004176D8" w:rsidRDefault="004176D8" w:rsidP="00E404D6">' +
    '<w:pPr><w:pStyle w:val="Lijstalinea"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Verhoog de &#8220;Current ATIS info&#8221; letter.</w:t></w:r>' +
    '<w:r><w:rPr><w:noProof/><w:lang w:eastAsia="nl-NL"/></w:rPr>' +
    '<w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0">' +
    '<wp:extent cx="5514975" cy="1885950"/><wp:effectExtent l="152400" t="152400" r="371475" b="361950"/>' +
    '<wp:docPr id="8" name="Afbeelding 8"/>' +
    '<wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr>' +
    '<a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture">' +
    '<pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="0" name="Picture 1"/>' +
    '<pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr>' +
    '<pic:blipFill rotWithShape="1"><a:blip r:embed="rId18"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}">' +
    '<a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip>' +
    '<a:srcRect t="69014"/><a:stretch/></pic:blipFill>' +
    '<pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="5514975" cy="1885950"/></a:xfrm>' +
    '<a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:ln><a:noFill/></a:ln>' +
    '<a:effectLst><a:outerShdw blurRad="292100" dist="139700" dir="2700000" algn="tl" rotWithShape="0">' +
    '<a:srgbClr val="333333"><a:alpha val="65000"/></a:srgbClr></a:outerShdw></a:effectLst>' +
    '<a:extLst><a:ext uri="{53640926-AAD7-44D8-BBD7-CCE9431645EC}">' +
    '<a14:shadowObscured xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main"/></a:ext></a:extLst>' +
    '</pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r>' +
    '</w:p>'
Set-ParagraphInnerXml $p6 $p6xml

Write-Host "Edits applied."
